# Update cryptocurrency price/volume data per Oct 23 2023 GitHub Actions refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.574.24"
$ws.Range("D3").Value = "1.669.09"
$ws.Range("E3").Value = "  +2.08%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "'218.77"
$ws.Range("E5").Value = "  +1.85%  "
$ws.Range("E6").Value = "  +1.75%  "
$ws.Range("D8").Value = "'29.10"
$ws.Range("E8").Value = "  +0.96%  "
$ws.Range("E9").Value = "  +1.96%  "
$ws.Range("D10").Value = "'0.0639"
$ws.Range("E10").Value = "  +4.85%  "
$ws.Range("E11").Value = "  -0.14%  "
$ws.Range("D12").Value = "1.911.32"
$ws.Range("E12").Value = "  +2.33%  "
$ws.Range("D13").Value = "1.670.93"
$ws.Range("E13").Value = "  +2.06%  "
$ws.Range("D14").Value = "'0.601"
$ws.Range("E14").Value = "  +6.60%  "
$ws.Range("E15").Value = "  +7.30%  "
$ws.Range("E16").Value = "  +3.94%  "
$ws.Range("D17").Value = "30.594.21"
$ws.Range("E17").Value = "  +2.11%  "
$ws.Range("D18").Value = "'65.87"
$ws.Range("E18").Value = "  +2.58%  "
$ws.Range("D19").Value = "'241.88"
$ws.Range("E19").Value = "  +0.05%  "
$ws.Range("E20").Value = "  +2.28%  "
$ws.Range("E21").Value = "  +0.17%  "
$ws.Range("E22").Value = "  +1.72%  "
$ws.Range("E23").Value = "  +1.10%  "
$ws.Range("E24").Value = "  -0.61%  "
$ws.Range("D25").Value = "'158.94"
$ws.Range("E25").Value = "  +0.66%  "
$ws.Range("D26").Value = "'15.76"
$ws.Range("E26").Value = "  +1.75%  "
$ws.Range("E27").Value = "  +2.15%  "
$ws.Range("E28").Value = "  +1.14%  "
$ws.Range("E29").Value = "  +0.11%  "
$ws.Range("E30").Value = "  +0.10%  "
$ws.Range("E31").Value = "  +3.40%  "
$ws.Range("D32").Value = "'3.45"
$ws.Range("E32").Value = "  +1.55%  "
$ws.Range("D33").Value = "'3.29"
$ws.Range("E33").Value = "  +3.43%  "
$ws.Range("D34").Value = "1.496.88"
$ws.Range("E34").Value = "  +4.70%  "
$ws.Range("E35").Value = "  +5.76%  "
$ws.Range("D36").Value = "'83.38"
$ws.Range("E36").Value = "  +10.24%  "
$ws.Range("D37").Value = "'1.01"
$ws.Range("E37").Value = "  -1.02%  "
$ws.Range("E38").Value = "  +7.66%  "
$ws.Range("E39").Value = "  +3.34%  "
$ws.Range("E40").Value = "  -3.82%  "
$ws.Range("D41").Value = "'2.29"
$ws.Range("E41").Value = "  +0.26%  "
$ws.Range("E42").Value = "  +0.75%  "
$ws.Range("E43").Value = "  +0.66%  "
$ws.Range("E44").Value = "  -0.17%  "
$ws.Range("B45").Value = "PaxDollar"
$ws.Range("C45").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D45").Value = "'1.00"
$ws.Range("E45").Value = "  +0.03%  "
$ws.Range("B46").Value = "WEMIXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D46").Value = "'1.01"
$ws.Range("E46").Value = "  +1.22%  "
$ws.Range("D47").Value = "'5.53"
$ws.Range("E47").Value = "  +3.40%  "
$ws.Range("D48").Value = "1.804.50"
$ws.Range("E48").Value = "  +1.68%  "
$ws.Range("D49").Value = "'49.64"
$ws.Range("E49").Value = "  -3.32%  "
$ws.Range("D50").Value = "'93.07"
$ws.Range("E50").Value = "  +2.72%  "
$ws.Range("E51").Value = "  +1.52%  "
